$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.442
$ws.Range("D7").Value = -7.453
$ws.Range("A10").Value = -21.846
$ws.Range("A12").Value = -21.589
$ws.Range("D15").Value = -8.175000000000001
$ws.Range("A18").Value = -22.167
$ws.Range("E18").Value = 16.208
$ws.Range("E19").Value = 16.538
$ws.Range("D20").Value = -7.57
$ws.Range("E27").Value = 16.257
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.236
$ws.Range("D31").Value = -8.096
$ws.Range("A37").Value = -20.21700000000001
$ws.Range("D40").Value = -7.93
$ws.Range("E42").Value = 16.665
$ws.Range("E44").Value = 16.62
$ws.Range("E47").Value = 16.467
$ws.Range("A55").Value = -22.283
$ws.Range("E58").Value = 16.387
$ws.Range("A68").Value = -21.534
$ws.Range("D68").Value = -6.778
$ws.Range("E73").Value = 16.635
$ws.Range("D76").Value = -7.753000000000002
$ws.Range("A77").Value = -20.963
$ws.Range("A78").Value = -20.281
$ws.Range("D87").Value = -8.257000000000001
$ws.Range("D88").Value = -8.295
$ws.Range("E95").Value = 17.397
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.052000000000001
$ws.Range("D101").Value = -7.616
$ws.Range("E101").Value = 16.44
$ws.Range("D102").Value = -8.142999999999999
